# Applies the cryptos-list price/volume refresh described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.182.41"
$ws.Range("E2").Value = "  -0.72%  "

$ws.Range("D3").Value = "1.826.72"
$ws.Range("E3").Value = "  -0.84%  "

$ws.Range("D4").Value = "'1.002"
$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").Value = "'234.00"

$ws.Range("D6").Value = "'0.5977"
$ws.Range("E6").Value = "  -4.37%  "

$ws.Range("D7").Value = "'1.003"
$ws.Range("E7").Value = "  +0.29%  "

$ws.Range("D8").Value = "'0.06957"
$ws.Range("E8").Value = "  -5.84%  "

$ws.Range("D9").Value = "'0.2748"
$ws.Range("E9").Value = "  -4.94%  "

$ws.Range("D10").Value = "'23.26"
$ws.Range("E10").Value = "  -6.49%  "

$ws.Range("D11").Value = "'0.07607"
$ws.Range("E11").Value = "  -1.40%  "

$ws.Range("D12").Value = "1.836.53"
$ws.Range("E12").Value = "  -0.25%  "

$ws.Range("E13").Value = "  -4.27%  "

$ws.Range("D14").Value = "'0.6252"
$ws.Range("E14").Value = "  -6.97%  "

$ws.Range("D15").Value = "'0.000009655"
$ws.Range("E15").Value = "  -5.82%  "

$ws.Range("D16").Value = "'78.19"
$ws.Range("E16").Value = "  -4.36%  "

$ws.Range("D17").Value = "28.864.64"
$ws.Range("E17").Value = "  -1.62%  "

$ws.Range("D18").Value = "'5.705"
$ws.Range("E18").Value = "  -9.18%  "

$ws.Range("D19").Value = "'220.84"
$ws.Range("E19").Value = "  -5.83%  "

$ws.Range("D20").Value = "'1.003"
$ws.Range("E20").Value = "  +0.24%  "

$ws.Range("D21").Value = "'11.52"
$ws.Range("E21").Value = "  -6.34%  "

$ws.Range("D22").Value = "'6.860"
$ws.Range("E22").Value = "  -5.98%  "

$ws.Range("D23").Value = "'1.004"
$ws.Range("E23").Value = "  +0.20%  "

$ws.Range("D24").Value = "'155.21"
$ws.Range("E24").Value = "  -1.24%  "

$ws.Range("D25").Value = "'7.947"
$ws.Range("E25").Value = "  -6.25%  "

$ws.Range("D26").Value = "'0.1285"
$ws.Range("E26").Value = "  -4.19%  "

$ws.Range("D27").Value = "'16.49"
$ws.Range("E27").Value = "  -4.74%  "

$ws.Range("D28").Value = "'0.06558"
$ws.Range("E28").Value = "  -10.10%  "

$ws.Range("D29").Value = "'1.455"
$ws.Range("E29").Value = "  -2.87%  "

$ws.Range("E30").Value = "  -2.41%  "

$ws.Range("D31").Value = "'3.835"
$ws.Range("E31").Value = "  -4.84%  "

$ws.Range("D32").Value = "'3.745"
$ws.Range("E32").Value = "  -7.37%  "

$ws.Range("D33").Value = "'1.091"
$ws.Range("E33").Value = "  -5.75%  "

$ws.Range("D34").Value = "'1.717"
$ws.Range("E34").Value = "  -5.56%  "

$ws.Range("D35").Value = "'0.6447"
$ws.Range("E35").Value = "  -9.29%  "

$ws.Range("D36").Value = "'2.539"
$ws.Range("E36").Value = "  -1.62%  "

$ws.Range("D37").Value = "'2.729"
$ws.Range("E37").Value = "  -1.99%  "

$ws.Range("D38").Value = "'0.01733"

$ws.Range("D39").Value = "'6.513"
$ws.Range("E39").Value = "  -4.04%  "

$ws.Range("D40").Value = "1.170.79"
$ws.Range("E40").Value = "  -5.13%  "

$ws.Range("D41").Value = "'0.8931"
$ws.Range("E41").Value = "  -6.30%  "

$ws.Range("D42").Value = "'1.002"
$ws.Range("E42").Value = "  +0.16%  "

$ws.Range("D43").Value = "1.979.49"
$ws.Range("E43").Value = "  -0.56%  "

$ws.Range("D44").Value = "'100.49"
$ws.Range("E44").Value = "  -0.63%  "

$ws.Range("D45").Value = "'61.93"
$ws.Range("E45").Value = "  -5.19%  "

$ws.Range("E46").Value = "  -3.53%  "

$ws.Range("B47").Value = "RenderToken"
$ws.Range("C47").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D47").Value = "'1.584"
$ws.Range("E47").Value = "  -6.77%  "

$ws.Range("D48").Value = "'8.432"
$ws.Range("E48").Value = "  -4.93%  "

$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.05534"
$ws.Range("E49").Value = "  -2.21%  "

$ws.Range("D50").Value = "'0.4544"
$ws.Range("E50").Value = "  -0.69%  "

$ws.Range("D51").Value = "'0.3641"
$ws.Range("E51").Value = "  -6.21%  "
